$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 792.1053000000001
$ws.Range("I80").Value = 1071.7142
$ws.Range("J80").Value = 629
$ws.Range("K80").Value = 3215.1426
$ws.Range("L80").Value = 1887
$ws.Range("M80").Value = -2217.1426
$ws.Range("N80").Value = -3883
$ws.Range("H83").Value = 792.1053000000001
$ws.Range("I83").Value = 1071.7142
$ws.Range("J83").Value = 629
$ws.Range("K83").Value = 9645.427799999999
$ws.Range("L83").Value = 5661
$ws.Range("M83").Value = -4653.427799999999
$ws.Range("N83").Value = -15645
$ws.Range("H92").Value = 649.9091
$ws.Range("I92").Value = 391.07407
$ws.Range("J92").Value = 1814.6666
$ws.Range("K92").Value = 391.07407
$ws.Range("L92").Value = 1814.6666
$ws.Range("M92").Value = 856.92593
$ws.Range("N92").Value = -4310.6666
$ws.Range("H100").Value = 1861.9
$ws.Range("I100").Value = 1266.6666
$ws.Range("J100").Value = 2754.75
$ws.Range("K100").Value = 1266.6666
$ws.Range("L100").Value = 2754.75
$ws.Range("M100").Value = -725.6666
$ws.Range("N100").Value = -3836.75
$ws.Range("H111").Value = 597.75
$ws.Range("I111").Value = 597.75
$ws.Range("K111").Value = 1793.25
$ws.Range("M111").Value = 1273.75
$ws.Range("H112").Value = 1607.4333
$ws.Range("I112").Value = 1043.25
$ws.Range("J112").Value = 1694.2307
$ws.Range("K112").Value = 3129.75
$ws.Range("L112").Value = 5082.6921
$ws.Range("M112").Value = -2021.75
$ws.Range("N112").Value = -7298.6921
$ws.Range("H113").Value = 22223882
$ws.Range("I113").Value = 28573144
$ws.Range("J113").Value = 1468
$ws.Range("K113").Value = 28573144
$ws.Range("L113").Value = 1468
$ws.Range("M113").Value = -28569890
$ws.Range("N113").Value = -7976
$ws.Range("H137").Value = 2343.2744
$ws.Range("I137").Value = 1968.84
$ws.Range("J137").Value = 2703.3076
$ws.Range("K137").Value = 5906.52
$ws.Range("L137").Value = 8109.9228
$ws.Range("M137").Value = -3356.52
$ws.Range("N137").Value = -13209.9228
$ws.Range("H141").Value = 1026.6666
$ws.Range("I141").Value = 1026.6666
$ws.Range("K141").Value = 3079.9998
$ws.Range("M141").Value = 2100.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1018.375
$ws.Range("I2").Value = 1018.375
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1018.375
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -905.375
$ws.Range("H61").Value = 90910350
$ws.Range("I61").Value = 111111980
$ws.Range("K61").Value = 111111980
$ws.Range("M61").Value = -111111768
$ws.Range("H110").Value = 566.6667
$ws.Range("I110").Value = 650
$ws.Range("K110").Value = 650
$ws.Range("M110").Value = 1395
$ws.Range("H116").Value = 1018.375
$ws.Range("I116").Value = 1018.375
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1018.375
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1275.625
$ws.Range("H132").Value = 2737.5957
$ws.Range("I132").Value = 2560.8
$ws.Range("J132").Value = 3253.25
$ws.Range("K132").Value = 7682.400000000001
$ws.Range("L132").Value = 9759.75
$ws.Range("M132").Value = -5152.400000000001
$ws.Range("N132").Value = -14819.75
$ws.Range("H136").Value = 90910350
$ws.Range("I136").Value = 111111980
$ws.Range("K136").Value = 333335940
$ws.Range("M136").Value = -333333390
$ws.Range("N2").ClearContents()
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1018.375
$ws.Range("I3").Value = 1018.375
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1018.375
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -904.375
$ws.Range("H86").Value = 4050.6
$ws.Range("I86").Value = 4201.231
$ws.Range("K86").Value = 4201.231
$ws.Range("M86").Value = -3078.231
$ws.Range("H89").Value = 4050.6
$ws.Range("I89").Value = 4201.231
$ws.Range("K89").Value = 21006.155
$ws.Range("M89").Value = -15390.155
$ws.Range("N3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 125001420
$ws.Range("I16").Value = 166668160
$ws.Range("J16").Value = 1206.5
$ws.Range("K16").Value = 166668160
$ws.Range("L16").Value = 1206.5
$ws.Range("M16").Value = -166667873
$ws.Range("N16").Value = -1780.5
$ws.Range("H31").Value = 1478.3469
$ws.Range("J31").Value = 1942.5
$ws.Range("L31").Value = 1942.5
$ws.Range("N31").Value = -2532.5
$ws.Range("H34").Value = 1478.3469
$ws.Range("J34").Value = 1942.5
$ws.Range("L34").Value = 1942.5
$ws.Range("N34").Value = -2346.5
$ws.Range("H58").Value = 7315.3887
$ws.Range("I58").Value = 1475.2307
$ws.Range("J58").Value = 22499.8
$ws.Range("K58").Value = 1475.2307
$ws.Range("L58").Value = 22499.8
$ws.Range("M58").Value = -1272.2307
$ws.Range("N58").Value = -22905.8
$ws.Range("H113").Value = 125001420
$ws.Range("I113").Value = 166668160
$ws.Range("J113").Value = 1206.5
$ws.Range("K113").Value = 166668160
$ws.Range("L113").Value = 1206.5
$ws.Range("M113").Value = -166665990
$ws.Range("N113").Value = -5546.5
$ws.Range("H132").Value = 1968.909
$ws.Range("I132").Value = 2146.1428
$ws.Range("J132").Value = 1658.75
$ws.Range("K132").Value = 6438.428400000001
$ws.Range("L132").Value = 4976.25
$ws.Range("M132").Value = -3908.428400000001
$ws.Range("N132").Value = -10036.25
$ws.Range("H134").Value = 14287029
$ws.Range("I134").Value = 1312.52
$ws.Range("J134").Value = 50001320
$ws.Range("K134").Value = 3937.56
$ws.Range("L134").Value = 150003960
$ws.Range("M134").Value = -1402.56
$ws.Range("N134").Value = -150009030
$ws.Range("H136").Value = 7315.3887
$ws.Range("I136").Value = 1475.2307
$ws.Range("J136").Value = 22499.8
$ws.Range("K136").Value = 4425.6921
$ws.Range("L136").Value = 67499.39999999999
$ws.Range("M136").Value = -1875.6921
$ws.Range("N136").Value = -72599.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 914.3
$ws.Range("I18").Value = 1051.6666
$ws.Range("J18").Value = 708.25
$ws.Range("K18").Value = 3154.9998
$ws.Range("L18").Value = 2124.75
$ws.Range("M18").Value = -2985.9998
$ws.Range("N18").Value = -2462.75
$ws.Range("H93").Value = 6754.8887
$ws.Range("J93").Value = 6754.8887
$ws.Range("L93").Value = 20264.6661
$ws.Range("N93").Value = -24008.6661
$ws.Range("H113").Value = 724.4524
$ws.Range("I113").Value = 664.1
$ws.Range("J113").Value = 743.3125
$ws.Range("K113").Value = 1992.3
$ws.Range("L113").Value = 2229.9375
$ws.Range("M113").Value = 177.6999999999998
$ws.Range("N113").Value = -6569.9375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 924.4737
$ws.Range("I122").Value = 925.8889
$ws.Range("J122").Value = 899
$ws.Range("K122").Value = 2777.6667
$ws.Range("L122").Value = 2697
$ws.Range("M122").Value = -327.6667000000002
$ws.Range("N122").Value = -7597
$ws.Range("H132").Value = 2722.0588
$ws.Range("I132").Value = 2272.5
$ws.Range("J132").Value = 3801
$ws.Range("K132").Value = 6817.5
$ws.Range("L132").Value = 11403
$ws.Range("M132").Value = -4287.5
$ws.Range("N132").Value = -16463

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2356.4285
$ws.Range("I7").Value = 2146.6667
$ws.Range("J7").Value = 2513.75
$ws.Range("K7").Value = 2146.6667
$ws.Range("L7").Value = 2513.75
$ws.Range("M7").Value = -2034.6667
$ws.Range("N7").Value = -2737.75
$ws.Range("H61").Value = 1254.4546
$ws.Range("I61").Value = 1157.7
$ws.Range("J61").Value = 2222
$ws.Range("K61").Value = 1157.7
$ws.Range("L61").Value = 2222
$ws.Range("M61").Value = -955.7
$ws.Range("N61").Value = -2626
$ws.Range("H82").Value = 2352.2222
$ws.Range("I82").Value = 2328.3333
$ws.Range("J82").Value = 2400
$ws.Range("K82").Value = 2328.3333
$ws.Range("L82").Value = 2400
$ws.Range("M82").Value = -1967.3333
$ws.Range("N82").Value = -3122
$ws.Range("H85").Value = 2352.2222
$ws.Range("I85").Value = 2328.3333
$ws.Range("J85").Value = 2400
$ws.Range("K85").Value = 2328.3333
$ws.Range("L85").Value = 2400
$ws.Range("M85").Value = -1080.3333
$ws.Range("N85").Value = -4896
$ws.Range("H113").Value = 1254.4546
$ws.Range("I113").Value = 1157.7
$ws.Range("J113").Value = 2222
$ws.Range("K113").Value = 1157.7
$ws.Range("L113").Value = 2222
$ws.Range("M113").Value = 1012.3
$ws.Range("N113").Value = -6562
$ws.Range("H126").Value = 2356.4285
$ws.Range("I126").Value = 2146.6667
$ws.Range("J126").Value = 2513.75
$ws.Range("K126").Value = 6440.000100000001
$ws.Range("L126").Value = 7541.25
$ws.Range("M126").Value = -3970.000100000001
$ws.Range("N126").Value = -12481.25
$ws.Range("H132").Value = 2639.2285
$ws.Range("I132").Value = 2516.5881
$ws.Range("K132").Value = 7549.7643
$ws.Range("M132").Value = -5019.7643

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1915.5927
$ws.Range("I132").Value = 1866.5
$ws.Range("K132").Value = 5599.5
$ws.Range("M132").Value = -3069.5
$ws.Range("H136").Value = 1854.4667
$ws.Range("I136").Value = 1661.2
$ws.Range("K136").Value = 4983.6
$ws.Range("M136").Value = -2433.6
